$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column G (header "K", formerly "Strike#") with recalculated values
$ws.Range("G2").Value = 6
$ws.Range("G3").Value = 6
$ws.Range("G4").Value = 7
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 3
